$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 10000.167
$ws.Range("J48").Value = 10000.167
$ws.Range("L48").Value = 30000.501
$ws.Range("N48").Value = -30584.501

$ws.Range("H56").Value = 10000.167
$ws.Range("J56").Value = 10000.167
$ws.Range("L56").Value = 30000.501
$ws.Range("N56").Value = -31068.501

$ws.Range("H132").Value = 899.4
$ws.Range("I132").Value = 888.6429000000001
$ws.Range("K132").Value = 2665.9287
$ws.Range("M132").Value = -135.9287000000004

$ws.Range("H137").Value = 7866.1904
$ws.Range("J137").Value = 8505.933999999999
$ws.Range("L137").Value = 25517.802
$ws.Range("N137").Value = -30617.802

$ws.Range("H138").Value = 1355413
$ws.Range("I138").Value = 2355.1428
$ws.Range("J138").Value = 2179013.5
$ws.Range("K138").Value = 7065.428400000001
$ws.Range("L138").Value = 6537040.5
$ws.Range("M138").Value = -1925.428400000001
$ws.Range("N138").Value = -6547320.5

$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3084427
$ws.Range("I32").Value = 3339046
$ws.Range("K32").Value = 3339046
$ws.Range("M32").Value = -3338759

$ws.Range("H61").Value = 37042668
$ws.Range("I61").Value = 1737.9474
$ws.Range("K61").Value = 1737.9474
$ws.Range("M61").Value = -1525.9474

$ws.Range("H122").Value = 5204.3335
$ws.Range("I122").Value = 3066.3333
$ws.Range("K122").Value = 9198.999899999999
$ws.Range("M122").Value = -6748.999899999999

$ws.Range("H132").Value = 5705.6216
$ws.Range("I132").Value = 1981.5714
$ws.Range("K132").Value = 5944.7142
$ws.Range("M132").Value = -3414.7142

$ws.Range("H136").Value = 37042668
$ws.Range("I136").Value = 1737.9474
$ws.Range("K136").Value = 5213.8422
$ws.Range("M136").Value = -2663.8422

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6671556
$ws.Range("I20").Value = 7249734.5
$ws.Range("J20").Value = 22499.5
$ws.Range("K20").Value = 7249734.5
$ws.Range("L20").Value = 22499.5
$ws.Range("M20").Value = -7249487.5
$ws.Range("N20").Value = -22993.5

$ws.Range("H86").Value = 8099031.5
$ws.Range("I86").Value = 11954305
$ws.Range("J86").Value = 2957
$ws.Range("K86").Value = 11954305
$ws.Range("L86").Value = 2957
$ws.Range("M86").Value = -11953182
$ws.Range("N86").Value = -5203

$ws.Range("H89").Value = 8099031.5
$ws.Range("I89").Value = 11954305
$ws.Range("J89").Value = 2957
$ws.Range("K89").Value = 59771525
$ws.Range("L89").Value = 14785
$ws.Range("M89").Value = -59765909
$ws.Range("N89").Value = -26017

$ws.Range("H134").Value = 15635593
$ws.Range("J134").Value = 11299.667
$ws.Range("L134").Value = 33899.001
$ws.Range("N134").Value = -38969.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3588.8696
$ws.Range("I16").Value = 868.73334
$ws.Range("J16").Value = 8689.125
$ws.Range("K16").Value = 868.73334
$ws.Range("L16").Value = 8689.125
$ws.Range("M16").Value = -581.73334
$ws.Range("N16").Value = -9263.125

$ws.Range("H31").Value = 6452.085
$ws.Range("I31").Value = 1990.1111
$ws.Range("K31").Value = 1990.1111
$ws.Range("M31").Value = -1695.1111

$ws.Range("H34").Value = 6452.085
$ws.Range("I34").Value = 1990.1111
$ws.Range("K34").Value = 1990.1111
$ws.Range("M34").Value = -1788.1111

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents() | Out-Null

$ws.Range("H58").Value = 6232.1953
$ws.Range("J58").Value = 9649
$ws.Range("L58").Value = 9649
$ws.Range("N58").Value = -10055

$ws.Range("H113").Value = 3588.8696
$ws.Range("I113").Value = 868.73334
$ws.Range("J113").Value = 8689.125
$ws.Range("K113").Value = 868.73334
$ws.Range("L113").Value = 8689.125
$ws.Range("M113").Value = 1301.26666
$ws.Range("N113").Value = -13029.125

$ws.Range("H132").Value = 4903.0557
$ws.Range("I132").Value = 2758.353
$ws.Range("K132").Value = 8275.059000000001
$ws.Range("M132").Value = -5745.059000000001

$ws.Range("H136").Value = 6232.1953
$ws.Range("J136").Value = 9649
$ws.Range("L136").Value = 28947
$ws.Range("N136").Value = -34047

$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 1746.138
$ws.Range("I15").Value = 13.625
$ws.Range("J15").Value = 2406.1428
$ws.Range("K15").Value = 40.875
$ws.Range("L15").Value = 7218.428400000001
$ws.Range("M15").Value = 99.125
$ws.Range("N15").Value = -7498.428400000001

$ws.Range("H107").Value = 16924246
$ws.Range("J107").Value = 17501266
$ws.Range("L107").Value = 52503798
$ws.Range("N107").Value = -52507638

$ws.Range("H113").Value = 4282.3076
$ws.Range("J113").Value = 5429.4736
$ws.Range("L113").Value = 16288.4208
$ws.Range("N113").Value = -20628.4208

$ws.Range("H128").Value = 115492.25
$ws.Range("I128").Value = 115492.25
$ws.Range("K128").Value = 346476.75
$ws.Range("M128").Value = -341496.75

$ws.Range("H132").Value = 13767.65
$ws.Range("I132").Value = 4195.1113
$ws.Range("J132").Value = 21599.727
$ws.Range("K132").Value = 37756.00169999999
$ws.Range("L132").Value = 194397.543
$ws.Range("M132").Value = -35226.00169999999
$ws.Range("N132").Value = -199457.543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents() | Out-Null

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents() | Out-Null

$ws.Range("H80").Value = 2500.7646
$ws.Range("I80").Value = 2246.818
$ws.Range("K80").Value = 2246.818
$ws.Range("M80").Value = -1248.818

$ws.Range("H83").Value = 2500.7646
$ws.Range("I83").Value = 2246.818
$ws.Range("K83").Value = 11234.09
$ws.Range("M83").Value = -6242.09

$ws.Range("H120").Value = 56134
$ws.Range("J120").Value = 56134
$ws.Range("L120").Value = 56134
$ws.Range("N120").Value = -65810

$ws.Range("H122").Value = 7265246.5
$ws.Range("J122").Value = 2002
$ws.Range("L122").Value = 6006
$ws.Range("N122").Value = -10906

$ws.Range("H126").Value = 5555.5415
$ws.Range("I126").Value = 2743.375
$ws.Range("J126").Value = 6961.625
$ws.Range("K126").Value = 8230.125
$ws.Range("L126").Value = 20884.875
$ws.Range("M126").Value = -5760.125
$ws.Range("N126").Value = -25824.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents() | Out-Null

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents() | Out-Null

$ws.Range("H94").Value = 57330
$ws.Range("J94").Value = 57330
$ws.Range("L94").Value = 57330
$ws.Range("N94").Value = -58682

$ws.Range("H100").Value = 2970.75
$ws.Range("I100").Value = 3021.2727
$ws.Range("J100").Value = 2859.6
$ws.Range("K100").Value = 3021.2727
$ws.Range("L100").Value = 2859.6
$ws.Range("M100").Value = -2480.2727
$ws.Range("N100").Value = -3941.6

$ws.Range("H136").Value = 11787.631
$ws.Range("I136").Value = 2912.7827
$ws.Range("J136").Value = 20662.479
$ws.Range("K136").Value = 8738.348100000001
$ws.Range("L136").Value = 61987.437
$ws.Range("M136").Value = -6188.348100000001
$ws.Range("N136").Value = -67087.43700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 100049
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 100049
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 100049
$ws.Range("M39").ClearContents() | Out-Null
$ws.Range("N39").Value = -100875

$ws.Range("H122").Value = 217274.73
$ws.Range("I122").Value = 404424
$ws.Range("K122").Value = 1213272
$ws.Range("M122").Value = -1210822

$ws.Range("H136").Value = 63132932
$ws.Range("I136").Value = 500002460
$ws.Range("J136").Value = 723000.4
$ws.Range("K136").Value = 1500007380
$ws.Range("L136").Value = 2169001.2
$ws.Range("M136").Value = -1500004830
$ws.Range("N136").Value = -2174101.2
